# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# with the latest scraped values (GitHub Actions refresh job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.309.67"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.48"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.18"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.76"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08030"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.53"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.940.20"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.888"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.70"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001030"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.323.21"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.526"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.152.91"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.32"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.039"
$ws.Range("E29").Value = "  +11.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.106"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.64"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.070"
$ws.Range("E32").Value = "  +7.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09498"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.426"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.543"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.384"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06079"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02245"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.229"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.172"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5867"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.513"
$ws.Range("E42").Value = "  +10.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1837"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07863"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.275"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5528"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.05"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.921"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.23"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.29"
$ws.Range("E51").Value = "  -2.63%  "
